$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) values from the crypto feed refresh.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "68.742.39"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.862.91"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.04%  "

$ws.Range("E4").Value = "  -0.05%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "600.24"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "162.22"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.861.41"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +3.03%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -1.69%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.168"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.30"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.17%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.459"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "36.85"
$cell.NumberFormat = "General"
$cell.Style = "Normal"

$ws.Range("E14").Value = "  -1.86%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.506.08"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.04%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.878.00"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +3.50%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "68.921.25"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "7.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.98%  "

$ws.Range("E19").Value = "  -0.46%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.13"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "11.34"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.60%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "483.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.0000162"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +6.77%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "83.88"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("E26").Value = "  -2.88%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "12.08"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("E28").Value = "  -0.07%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.93"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("E30").Value = "  -0.99%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.014.51"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "

$ws.Range("E32").Value = "  -3.39%  "

$ws.Range("E33").Value = "  +2.41%  "

$ws.Range("E34").Value = "  -4.22%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.806.66"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +3.42%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.03"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  -1.32%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("E42").Value = "  -2.44%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "430.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "48.48"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("E47").Value = "  -1.02%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "143.23"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.837.58"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.77%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0357"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "25.91"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +13.19%  "
